$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New feature request: User Profiles
$ws.Range("A10").Value = "User Profiles"

# Row 5 ("Roll groupings") gains two requestors that were missing before.
$ws.Range("D5").Value = "_dharwin - Reddit /r/dnd"
$ws.Range("E5").Value = "BlankTheorist - Reddit /r/dnd"

$ws.Range("B10").Value = "Allow for saving custom rolls into a user profile so that they do not bleed into other rolls"
$ws.Range("D10").Value = "_dharwin - Reddit /r/dnd"

# New feature request: Roll Negative Die
$ws.Range("A11").Value = "Roll Negative Die"
$ws.Range("B11").Value = "In order to roll with Bane, I need to roll a negative d4. Right now you can only roll positively"
$ws.Range("D11").Value = "Weston Fiala"

# Match the author's final cursor position recorded in the saved file.
$ws.Range("B13").Select()
